# Apply updates described by the commit diff to cryptos.xlsx
# The sheet stores every data value (price, volume%, coin name, link) as
# plain text in inline strings, even though columns D and E look numeric.
# To avoid Excel auto-converting these into real numbers/percentages
# (which would silently change the underlying type and drop things like
# trailing zeros, e.g. "0.05620" -> 0.0562), we force the target cells to
# a literal "Text" number format ("@") before writing the string value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) and column E (Volume(1h)) updates - these are numeric-
# or percent-looking strings, so they must stay as literal text.
$priceAndVolumeUpdates = @(
    @("D2", "244.03"),
    @("E2", "-0.61%"),
    @("D3", "27.13"),
    @("E3", "3.95%"),
    @("D4", "5.155"),
    @("E4", "0.58%"),
    @("D5", "0.05620"),
    @("E5", "0.58%"),
    @("D6", "6.481"),
    @("E6", "-0.20%"),
    @("D7", "0.8162"),
    @("E7", "0.05%"),
    @("D8", "0.8325"),
    @("E8", "-1.78%"),
    @("D9", "0.1330"),
    @("E9", "-0.82%"),
    @("D10", "0.06910"),
    @("E10", "-0.68%"),
    @("D11", "0.02888"),
    @("E11", "1.28%"),
    @("D12", "0.09383"),
    @("E12", "-0.15%"),
    @("D13", "0.001527"),
    @("E13", "0.06%"),
    @("D14", "0.0005968"),
    @("E14", "-0.72%"),
    @("D15", "0.006170"),
    @("E15", "-1.68%"),
    @("E16", "1.67%"),
    @("D17", "3.022"),
    @("E17", "-0.19%"),
    @("D18", "2.309"),
    @("E18", "8.99%"),
    @("E19", "-0.71%"),
    @("D20", "0.03088"),
    @("E20", "-4.05%"),
    @("E21", "-2.14%"),
    @("D22", "3.738"),
    @("E22", "-0.07%"),
    @("D23", "0.04494"),
    @("E23", "-4.65%"),
    @("E24", "-0.10%"),
    @("D25", "0.001223"),
    @("E25", "-2.04%"),
    @("D26", "0.004488"),
    @("E26", "-2.47%"),
    @("D27", "0.00009794"),
    @("E27", "1.99%"),
    @("E28", "0.62%"),
    @("E40", "-0.53%"),
    @("D41", "0.006056"),
    @("E41", "-2.97%"),
    @("D42", "0.1050"),
    @("E42", "-0.36%"),
    @("D43", "0.002608"),
    @("E43", "4.99%"),
    @("D44", "0.008199"),
    @("E44", "4.01%"),
    @("D45", "0.00005316"),
    @("E45", "0.11%"),
    @("E46", "-0.09%"),
    @("D47", "0.1090"),
    @("E47", "-18.39%"),
    @("D48", "0.002491"),
    @("E48", "21.60%"),
    @("E49", "-0.09%"),
    @("E50", "-0.09%")
)

foreach ($item in $priceAndVolumeUpdates) {
    $cellRef = $item[0]
    $newValue = $item[1]
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $newValue
}

# Column B (Coin) and column C (Link) updates - plain text, no special
# number-format handling required.
$coinAndLinkUpdates = @(
    @("B41", "KickToken"),
    @("C41", "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"),
    @("B42", "BKEXToken"),
    @("C42", "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"),
    @("B43", "CEJI"),
    @("C43", "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji")
)

foreach ($item in $coinAndLinkUpdates) {
    $cellRef = $item[0]
    $newValue = $item[1]
    $ws.Range($cellRef).Value = $newValue
}
